$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme
try { Write-Output ("Name=" + $tcs.Name) } catch { Write-Output ("ERR: " + $_.Exception.Message) }
try { $tcs.Name = "Office"; Write-Output ("Name2=" + $tcs.Name) } catch { Write-Output ("ERR2: " + $_.Exception.Message) }
